$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '51.580.83'
$ws.Range('E2').Value = '  +0.86%  '
$ws.Range('D3').Value = '2.983.35'
$ws.Range('E3').Value = '  +2.42%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.12%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '380.40'
$ws.Range('E5').Value = '  +4.26%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '106.24'
$ws.Range('E6').Value = '  +2.68%  '
$ws.Range('E7').Value = '  +0.87%  '
$ws.Range('E8').Value = '  -0.08%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.599'
$ws.Range('E9').Value = '  +1.46%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '37.45'
$ws.Range('E10').Value = '  +1.25%  '
$ws.Range('E11').Value = '  +0.43%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0844'
$ws.Range('E12').Value = '  +1.09%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '18.69'
$ws.Range('E13').Value = '  +1.28%  '
$ws.Range('D14').Value = '3.451.52'
$ws.Range('E14').Value = '  +2.21%  '
$ws.Range('E15').Value = '  +2.12%  '
$ws.Range('D16').Value = '2.980.31'
$ws.Range('E16').Value = '  +2.21%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.970'
$ws.Range('E17').Value = '  +2.00%  '
$ws.Range('D18').Value = '51.585.99'
$ws.Range('E18').Value = '  +0.86%  '
$ws.Range('E19').Value = '  +2.39%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.42'
$ws.Range('E20').Value = '  +2.30%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.03'
$ws.Range('E21').Value = '  -0.02%  '
$ws.Range('D22').Value = '0.0₃0960'
$ws.Range('E22').Value = '  +1.44%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '69.28'
$ws.Range('E23').Value = '  +1.76%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '263.92'
$ws.Range('E24').Value = '  +1.55%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.80'
$ws.Range('E25').Value = '  +4.05%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.172'
$ws.Range('E26').Value = '  -2.12%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.23'
$ws.Range('E27').Value = '  +18.04%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.44'
$ws.Range('E28').Value = '  +2.10%  '
$ws.Range('E29').Value = '  +0.02%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '26.07'
$ws.Range('E30').Value = '  +0.47%  '
$ws.Range('E31').Value = '  +3.32%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '9.90'
$ws.Range('E32').Value = '  -0.44%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '34.78'
$ws.Range('E33').Value = '  -1.17%  '
$ws.Range('B34').Value = 'VeChain'
$ws.Range('C34').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0465'
$ws.Range('E34').Value = '  +10.18%  '
$ws.Range('B35').Value = 'Toncoin'
$ws.Range('C35').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.10'
$ws.Range('E35').Value = '  -2.16%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '51.43'
$ws.Range('E36').Value = '  +1.47%  '
$ws.Range('E37').Value = '  -0.03%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.10'
$ws.Range('E38').Value = '  -1.26%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '17.44'
$ws.Range('E39').Value = '  +2.77%  '
$ws.Range('E40').Value = '  -6.83%  '
$ws.Range('E41').Value = '  -0.64%  '
$ws.Range('E42').Value = '  +2.50%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '123.72'
$ws.Range('E43').Value = '  +4.95%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '22.19'
$ws.Range('E44').Value = '  -1.59%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.08'
$ws.Range('E45').Value = '  -0.90%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.280'
$ws.Range('E46').Value = '  +18.85%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.38'
$ws.Range('E47').Value = '  +5.01%  '
$ws.Range('D48').Value = '2.049.81'
$ws.Range('E48').Value = '  -0.88%  '
$ws.Range('E49').Value = '  +1.59%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0353'
$ws.Range('E50').Value = '  +10.28%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '5.19'
$ws.Range('E51').Value = '  +3.36%  '
